$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was captured for "Ají" (Inferno variety) on
# Región de Arica y Parinacota; insert it as a new row 8, pushing the
# existing rows 8-97 down to 9-98 (dimension grows from A1:R97 to A1:R98).
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8, 3).Value = "Ñuble"
$ws.Cells.Item(8, 4).Value = 44817
$ws.Cells.Item(8, 5).Value = 16
$ws.Cells.Item(8, 6).Value = 100112021
$ws.Cells.Item(8, 7).Value = "Ají"
$ws.Cells.Item(8, 8).Value = "Inferno"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 60
$ws.Cells.Item(8, 11).Value = 22000
$ws.Cells.Item(8, 12).Value = 23000
$ws.Cells.Item(8, 13).Value = 22500
$ws.Cells.Item(8, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 1500
$ws.Cells.Item(8, 17).Value = 15
$ws.Cells.Item(8, 18).Value = "Hortaliza"
